$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.481.38"
$ws.Range("E2").Value = "  +5.15%  "
$ws.Range("D3").Value = "2.057.03"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.05%  "
$ws.Range("E6").Value = "  +3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.59%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0770"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.935"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +25.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "2.361.75"
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.20%  "
$ws.Range("D18").Value = "2.067.44"
$ws.Range("E18").Value = "  +2.44%  "
$ws.Range("D19").Value = "37.440.82"
$ws.Range("E19").Value = "  +5.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "0.0₃0879"
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("E26").Value = "  +7.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +45.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.97%  "
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.01%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.53%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.24%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +38.27%  "
$ws.Range("E41").Value = "  +14.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.18%  "
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("E44").Value = "  +6.67%  "
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Value = "1.418.12"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("E51").Value = "  +10.10%  "
